$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("referente")
$ws2 = $wb.Worksheets.Item("obra")

# Implementacion de Dublin Core: rename Spanish metadata headers to
# their Dublin Core English equivalents.
$ws1.Range("B1").Value = "title"
$ws1.Range("C1").Value = "date"
$ws1.Range("D1").Value = "publisher"

$ws2.Range("B1").Value = "title"
$ws2.Range("C1").Value = "date"
$ws2.Range("D1").Value = "format"
$ws2.Range("E1").Value = "medium"

# Remove the now-redundant "bga-obra" sheet (merge of the two above).
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("bga-obra").Delete()
$excel.DisplayAlerts = $true

# Restore selections on both remaining sheets.
$ws2.Range("E10").Select()
$ws1.Activate()
$ws1.Range("B17").Select()
